$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 139.28572
$ws.Range("I6").Value = 80
$ws.Range("K6").Value = 240
$ws.Range("M6").Value = -128
$ws.Range("H19").Value = 672.6
$ws.Range("I19").Value = 301
$ws.Range("K19").Value = 301
$ws.Range("M19").Value = -126
$ws.Range("H28").Value = 2149.9
$ws.Range("I28").Value = 1784.8667
$ws.Range("J28").Value = 3245
$ws.Range("K28").Value = 1784.8667
$ws.Range("L28").Value = 3245
$ws.Range("M28").Value = -1299.8667
$ws.Range("N28").Value = -4215
$ws.Range("H86").Value = 16670507
$ws.Range("I86").Value = 3414.1428
$ws.Range("K86").Value = 3414.1428
$ws.Range("M86").Value = -2291.1428
$ws.Range("H89").Value = 16670507
$ws.Range("I89").Value = 3414.1428
$ws.Range("K89").Value = 17070.714
$ws.Range("M89").Value = -11454.714
$ws.Range("H92").Value = 1202.9375
$ws.Range("I92").Value = 1110.9231
$ws.Range("K92").Value = 1110.9231
$ws.Range("M92").Value = 137.0769
$ws.Range("H100").Value = 2576
$ws.Range("I100").Value = 2000
$ws.Range("J100").Value = 3152
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 3152
$ws.Range("M100").Value = -1459
$ws.Range("N100").Value = -4234
$ws.Range("H135").Value = 4805.5557
$ws.Range("I135").Value = 250
$ws.Range("J135").Value = 5375
$ws.Range("K135").Value = 2250
$ws.Range("L135").Value = 48375
$ws.Range("M135").Value = 285
$ws.Range("N135").Value = -53445
$ws.Range("H137").Value = 3948.889
$ws.Range("I137").Value = 1390.1666
$ws.Range("J137").Value = 9066.333000000001
$ws.Range("K137").Value = 4170.4998
$ws.Range("L137").Value = 27198.999
$ws.Range("M137").Value = -1620.4998
$ws.Range("N137").Value = -32298.999
$ws.Range("H141").Value = 12195.77
$ws.Range("I141").Value = 12878.75
$ws.Range("K141").Value = 38636.25
$ws.Range("M141").Value = -33456.25

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14707462
$ws.Range("I32").Value = 16394584
$ws.Range("K32").Value = 16394584
$ws.Range("M32").Value = -16394297
$ws.Range("H69").Value = 156972.67
$ws.Range("J69").Value = 156972.67
$ws.Range("L69").Value = 156972.67
$ws.Range("N69").Value = -158470.67
$ws.Range("H72").Value = 156972.67
$ws.Range("J72").Value = 156972.67
$ws.Range("L72").Value = 470918.01
$ws.Range("N72").Value = -478406.01
$ws.Range("H74").Value = 100003880
$ws.Range("I74").Value = 166669870
$ws.Range("K74").Value = 166669870
$ws.Range("M74").Value = -166668996
$ws.Range("H77").Value = 100003880
$ws.Range("I77").Value = 166669870
$ws.Range("K77").Value = 833349350
$ws.Range("M77").Value = -833344982
$ws.Range("H102").Value = 1955.2858
$ws.Range("I102").Value = 1955.2858
$ws.Range("K102").Value = 1955.2858
$ws.Range("M102").Value = -333.2858000000001
$ws.Range("H122").Value = 8338178
$ws.Range("I122").Value = 3705.75
$ws.Range("J122").Value = 9264230
$ws.Range("K122").Value = 11117.25
$ws.Range("L122").Value = 27792690
$ws.Range("M122").Value = -8667.25
$ws.Range("N122").Value = -27797590

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 392.16666
$ws.Range("I11").Value = 369.57144
$ws.Range("J11").Value = 471.25
$ws.Range("K11").Value = 369.57144
$ws.Range("L11").Value = 471.25
$ws.Range("M11").Value = -229.57144
$ws.Range("N11").Value = -751.25
$ws.Range("H26").Value = 16425.637
$ws.Range("I26").Value = 12516.3
$ws.Range("K26").Value = 12516.3
$ws.Range("M26").Value = -12224.3
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H107").Value = 1904.4546
$ws.Range("I107").Value = 1844.9
$ws.Range("K107").Value = 1844.9
$ws.Range("M107").Value = 75.09999999999991
$ws.Range("H132").Value = 63972.223
$ws.Range("J132").Value = 63972.223
$ws.Range("L132").Value = 63972.223
$ws.Range("N132").Value = -74092.223

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H58").Value = 5879.6924
$ws.Range("I58").Value = 6193.6
$ws.Range("K58").Value = 6193.6
$ws.Range("M58").Value = -5990.6
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H132").Value = 1685.3513
$ws.Range("I132").Value = 1666.9429
$ws.Range("J132").Value = 2007.5
$ws.Range("K132").Value = 5000.8287
$ws.Range("L132").Value = 6022.5
$ws.Range("M132").Value = -2470.8287
$ws.Range("N132").Value = -11082.5
$ws.Range("H134").Value = 3232.5122
$ws.Range("I134").Value = 2622.6667
$ws.Range("J134").Value = 5748.125
$ws.Range("K134").Value = 7868.000100000001
$ws.Range("L134").Value = 17244.375
$ws.Range("M134").Value = -5333.000100000001
$ws.Range("N134").Value = -22314.375
$ws.Range("H136").Value = 5879.6924
$ws.Range("I136").Value = 6193.6
$ws.Range("K136").Value = 18580.8
$ws.Range("M136").Value = -16030.8

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 17009.572
$ws.Range("I24").Value = 10480.2
$ws.Range("J24").Value = 33333
$ws.Range("K24").Value = 10480.2
$ws.Range("L24").Value = 33333
$ws.Range("M24").Value = -10307.2
$ws.Range("N24").Value = -33679
$ws.Range("H113").Value = 4198.625
$ws.Range("J113").Value = 5000
$ws.Range("L113").Value = 5000
$ws.Range("N113").Value = -9340

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 419.2
$ws.Range("I55").Value = 385.2857
$ws.Range("K55").Value = 385.2857
$ws.Range("M55").Value = -212.2857
$ws.Range("H61").Value = 3697.875
$ws.Range("I61").Value = 2847.3333
$ws.Range("J61").Value = 6249.5
$ws.Range("K61").Value = 2847.3333
$ws.Range("L61").Value = 6249.5
$ws.Range("M61").Value = -2645.3333
$ws.Range("N61").Value = -6653.5
$ws.Range("H113").Value = 3697.875
$ws.Range("I113").Value = 2847.3333
$ws.Range("J113").Value = 6249.5
$ws.Range("K113").Value = 2847.3333
$ws.Range("L113").Value = 6249.5
$ws.Range("M113").Value = -677.3332999999998
$ws.Range("N113").Value = -10589.5
$ws.Range("H122").Value = 6949025
$ws.Range("I122").Value = 4312.615
$ws.Range("J122").Value = 25005278
$ws.Range("K122").Value = 12937.845
$ws.Range("L122").Value = 75015834
$ws.Range("M122").Value = -10487.845
$ws.Range("N122").Value = -75020734

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 4059000.5
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
